# Update the "取得日時" (retrieved-at) timestamp column for the existing
# data rows on the "ランサーズ" sheet from 2025-10-19 06:25:53 to
# 2025-10-19 06:32:17 (latest scrape run).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-10-19 06:25:53"
$newTimestamp = "2025-10-19 06:32:17"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
